# Applies the "Leitura de Titulo e Abstract / Terminado" edit:
#  - Fills in column F (Critério de exclusão = 5) and column G
#    (Status - Etapa 1 = "Eliminado") for the rows that were reviewed.
#  - Adjusts a few column widths / the saved sheet-view scroll position
#    to match the state the workbook was left in after the review pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only column F (exclusion-criterion code) needs to be set,
# because column G already had a "Eliminado" status.
$fOnlyRows = @(8, 16, 46)
foreach ($r in $fOnlyRows) {
    $ws.Cells.Item($r, 6).Value = 5
}

# Row 17 only gets the exclusion-criterion code, no status value.
$ws.Cells.Item(17, 6).Value = 5

# Rows where both column F and column G need to be populated.
$bothRows = @(5, 9, 11, 21, 23, 28, 30, 33, 37, 38, 39, 43, 44, 45, 47)
foreach ($r in $bothRows) {
    $ws.Cells.Item($r, 6).Value = 5
    $ws.Cells.Item($r, 7).Value = "Eliminado"
}

# Column width adjustments (columns C, D, F, H) - the hidden columns
# D, F and H were unhidden and resized, C was widened.
$ws.Columns.Item(3).ColumnWidth = 71.66
$ws.Columns.Item(4).ColumnWidth = 18.66
$ws.Columns.Item(4).Hidden = $false
$ws.Columns.Item(6).ColumnWidth = 11.83
$ws.Columns.Item(6).Hidden = $false
$ws.Columns.Item(8).ColumnWidth = 8.0
$ws.Columns.Item(8).Hidden = $false

# Restore the view to the top of the sheet and move the active selection.
$ws.Range("C1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G45").Select()
